# Updated data with DX explanations
# Adds a new "Participant Status" row (row 11) to the Basic_Demos data
# dictionary sheet, with a Values / Value Labels explanation (wrapped,
# taller row), and moves the active selection down past the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string pool order must come out as:
#   28 Participant_Status
#   29 Participant Status
#   30 Complete, \nDropout, \nLost to follow-up
#   31 Complete=...\nDropout=...\nLost to follow-up=...
# so populate column B (Variable Name) before column A (Question).
$ws.Range("B11").Value = "Participant_Status"
$ws.Range("A11").Value = "Participant Status"
$ws.Range("C11").Value = "Text"
$ws.Range("D11").Value = "Complete, `nDropout, `nLost to follow-up"
$ws.Range("E11").Value = "Complete=Participant completed study`nDropout=Participant dropped out of study`nLost to follow-up=Participant didn't follow up with study coordinators"

# D11/E11 carry the wrapped "Values" / "Value Labels" style used elsewhere
# in the sheet (e.g. D10/E10, E8).
$ws.Range("D11:E11").WrapText = $true

# Row grew tall enough to show all the wrapped label text.
$ws.Rows.Item(11).RowHeight = 96

# The author's selection ended up one row below the newly typed data.
$ws.Range("E12").Select()
